$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 90, shifting existing rows 90-209 down to 91-210
$ws.Rows(90).Insert()

# Populate the newly inserted row 90 with its data
$ws.Range("A90").Value = 11
$ws.Range("B90").Value = "Vega Monumental Concepción"
$ws.Range("C90").Value = "Bíobío"
$ws.Range("D90").Value = 44482
$ws.Range("E90").Value = 8
$ws.Range("F90").Value = 100114014
$ws.Range("G90").Value = "Betarraga"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 450
$ws.Range("K90").Value = 550
$ws.Range("L90").Value = 600
$ws.Range("M90").Value = 572
$ws.Range("N90").Value = "$/paquete 5 unidades"
$ws.Range("O90").Value = "Región del Maule"
$ws.Range("P90").Value = 114
$ws.Range("Q90").Value = 5
$ws.Range("R90").Value = "Hortaliza"
